$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = 2.34
$ws.Range("L3").Value = 1.39
$ws.Range("W3").Value = 1.74
$ws.Range("AF3").Value = 15.5
$ws.Range("AN3").Value = 17.5
$ws.Range("H4").Value = 5.5
$ws.Range("I4").Value = 5.6
$ws.Range("J4").Value = 3.85
$ws.Range("K4").Value = 3.9
$ws.Range("O4").Value = 1.36
$ws.Range("U4").Value = 1.93
$ws.Range("F5").Value = 1.84
$ws.Range("G5").Value = 1.86
$ws.Range("AN5").Value = 12
$ws.Range("G6").Value = 2.08
$ws.Range("H6").Value = 3.65
$ws.Range("I6").Value = 3.8
$ws.Range("J6").Value = 4
$ws.Range("P6").Value = 2.36
$ws.Range("Q6").Value = 1.62
$ws.Range("S6").Value = 2.52
$ws.Range("T6").Value = 1.61
$ws.Range("Y6").Value = 23
$ws.Range("AF6").Value = 18
$ws.Range("AG6").Value = 11
$ws.Range("AJ6").Value = 29
$ws.Range("AN6").Value = 11
$ws.Range("G8").Value = 2.06
$ws.Range("Q8").Value = 2.14
$ws.Range("Q9").Value = 2.08
$ws.Range("F10").Value = 2.54
$ws.Range("G10").Value = 3.1
$ws.Range("H10").Value = 2.54
$ws.Range("I10").Value = 3.1
$ws.Range("J10").Value = 3.55
$ws.Range("P10").Value = 2.2
$ws.Range("Q10").Value = 1.66
$ws.Range("Q11").Value = 2.7
$ws.Range("S12").Value = 4.1
$ws.Range("AF12").Value = 11
$ws.Range("G13").Value = 1.15
$ws.Range("J13").Value = 10.5
$ws.Range("F14").Value = 1.97
$ws.Range("G14").Value = 2.46
$ws.Range("K14").Value = 5.4
$ws.Range("P14").Value = 1.63
$ws.Range("G15").Value = 2.06
$ws.Range("H15").Value = 3.4
$ws.Range("F16").Value = 4.3
$ws.Range("H16").Value = 1.9
$ws.Range("I16").Value = 2
$ws.Range("J16").Value = 3.6
$ws.Range("N16").Value = 3.55
$ws.Range("O16").Value = 1.32
$ws.Range("P16").Value = 1.87
$ws.Range("AO16").Value = 17.5
$ws.Range("G18").Value = 1.98
$ws.Range("F19").Value = 2.56
$ws.Range("G19").Value = 2.74
$ws.Range("H19").Value = 2.84
$ws.Range("H20").Value = 7.2
$ws.Range("I20").Value = 7.6
$ws.Range("K20").Value = 4.8
$ws.Range("T20").Value = 1.91
$ws.Range("F21").Value = 2.6
$ws.Range("G21").Value = 2.64
$ws.Range("H21").Value = 2.98
$ws.Range("I21").Value = 3.05
$ws.Range("N21").Value = 3.75
$ws.Range("O21").Value = 1.34
$ws.Range("P21").Value = 1.91
$ws.Range("Q21").Value = 2.04
$ws.Range("R21").Value = 1.36
$ws.Range("S21").Value = 3.6
$ws.Range("T21").Value = 1.8
$ws.Range("U21").Value = 2.18
$ws.Range("X21").Value = 14
$ws.Range("Z21").Value = 21
$ws.Range("AA21").Value = 55
$ws.Range("AB21").Value = 11
$ws.Range("AD21").Value = 13.5
$ws.Range("AF21").Value = 17
$ws.Range("AG21").Value = 12.5
$ws.Range("AH21").Value = 17.5
$ws.Range("AI21").Value = 46
$ws.Range("AM21").Value = 85
$ws.Range("AN21").Value = 23
$ws.Range("F24").Value = 1.33
$ws.Range("G24").Value = 1.37
$ws.Range("H24").Value = 10
$ws.Range("I24").Value = 12.5
$ws.Range("J24").Value = 5.4
$ws.Range("K24").Value = 6.2
$ws.Range("P24").Value = 2.26
$ws.Range("Q24").Value = 1.67
$ws.Range("T24").Value = 2.08
$ws.Range("U24").Value = 1.8
$ws.Range("X24").Value = 1000
$ws.Range("AB24").Value = 1000
$ws.Range("AC24").Value = 1000
$ws.Range("AE24").Value = 230
$ws.Range("AF24").Value = 980
$ws.Range("AG24").Value = 11
$ws.Range("AI24").Value = 180
$ws.Range("AJ24").Value = 1000
$ws.Range("AK24").Value = 1000
$ws.Range("AM24").Value = 210
$ws.Range("AN24").Value = 5.6
$ws.Range("I25").Value = 5.8
$ws.Range("G27").Value = 2.88
$ws.Range("I28").Value = 3.9
$ws.Range("K28").Value = 4.2
$ws.Range("H30").Value = 6.2
$ws.Range("I30").Value = 17
$ws.Range("J30").Value = 4.7
$ws.Range("K30").Value = 6.8
$ws.Range("P30").Value = 2.24
$ws.Range("Q30").Value = 1.63
$ws.Range("I33").Value = 7.8
$ws.Range("H34").Value = 2.66
$ws.Range("I34").Value = 2.74
$ws.Range("S34").Value = 5.2
$ws.Range("F35").Value = 1.79
$ws.Range("K35").Value = 4.2
$ws.Range("P35").Value = 2.04
$ws.Range("Q35").Value = 1.81
$ws.Range("G36").Value = 1.64
$ws.Range("Q36").Value = 1.8
$ws.Range("G37").Value = 32
$ws.Range("H37").Value = 1.36
$ws.Range("P38").Value = 1.83
$ws.Range("U39").Value = 2.56
$ws.Range("F42").Value = 3.05
$ws.Range("G42").Value = 3.85
$ws.Range("H42").Value = 2.02
$ws.Range("I42").Value = 2.48
$ws.Range("J42").Value = 2.96
$ws.Range("K42").Value = 4.4
$ws.Range("P42").Value = 2.22
$ws.Range("Q42").Value = 1.64
$ws.Range("J43").Value = 4.6
$ws.Range("N44").Value = 4
$ws.Range("P44").Value = 2.06
$ws.Range("F47").Value = 1.78
$ws.Range("H47").Value = 4.3
$ws.Range("I47").Value = 5.2
$ws.Range("J47").Value = 3.7
$ws.Range("K47").Value = 4.4
$ws.Range("M47").Value = 1.01
$ws.Range("N47").Value = 2.04
$ws.Range("O47").Value = 1.28
$ws.Range("P47").Value = 2.04
$ws.Range("Q47").Value = 1.78
$ws.Range("R47").Value = 1.31
$ws.Range("S47").Value = 2.74
$ws.Range("T47").Value = 1.01
$ws.Range("U47").Value = 1.01
$ws.Range("X47").Value = 25
$ws.Range("Y47").Value = 26
$ws.Range("Z47").Value = 55
$ws.Range("AA47").Value = 1000
$ws.Range("AB47").Value = 14.5
$ws.Range("AC47").Value = 13
$ws.Range("AD47").Value = 28
$ws.Range("AE47").Value = 85
$ws.Range("AF47").Value = 17.5
$ws.Range("AG47").Value = 15
$ws.Range("AH47").Value = 28
$ws.Range("AI47").Value = 90
$ws.Range("AJ47").Value = 30
$ws.Range("AK47").Value = 28
$ws.Range("AL47").Value = 48
$ws.Range("AM47").Value = 1000
$ws.Range("AN47").Value = 1000
$ws.Range("AO47").Value = 1000
$ws.Range("F49").Value = 2.36
$ws.Range("H49").Value = 3.3
$ws.Range("I49").Value = 3.6
$ws.Range("F50").Value = 2.92
$ws.Range("H50").Value = 2.56
$ws.Range("I50").Value = 2.62
$ws.Range("J50").Value = 3.55
$ws.Range("K50").Value = 3.7
$ws.Range("H51").Value = 12.5
$ws.Range("I51").Value = 26
$ws.Range("H52").Value = 1.52
$ws.Range("N53").Value = 3.45
$ws.Range("U54").Value = 2.62
$ws.Range("F55").Value = 2.16
$ws.Range("G55").Value = 2.18
$ws.Range("K55").Value = 3.4
$ws.Range("P55").Value = 1.6
$ws.Range("H56").Value = 4.3
$ws.Range("J57").Value = 3.15
